$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For plotting")

$ws.Range("D2").Value = 0.0536867743516621
$ws.Range("E2").Value = 0.163641117886128

$ws.Range("D3").Value = -0.00826467306677648
$ws.Range("E3").Value = 0.115426697884501

$ws.Range("D4").Value = 0.0042012509749491
$ws.Range("E4").Value = 0.131851978652704
